$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Bugs and errors": only the cursor/selection moved (cosmetic).
# ---------------------------------------------------------------------------
$wsBugs = $wb.Worksheets.Item("Bugs and errors")
$wsBugs.Activate()
$wsBugs.Range("B26").Select()

# ---------------------------------------------------------------------------
# Sheet "Implemented Features": new "var/mean of inning" work log entries.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Implemented Features")
$ws.Activate()

# The "Commited on" / "Git branch" note that used to live on row 5 (E5/F5)
# now belongs to the newly added row 7 (E7/F7) - clear it from row 5 first.
$ws.Range("E5:F5").Clear()

# New row 6: "Added missing data catches in mains"
$ws.Range("A6").Value = "Added missing data catches in mains"
$ws.Range("B6").Value = 44645
$ws.Range("C6").Value = "Thomas"

# New row 7: "Var / mean of inning"
$ws.Range("A7").Value = "Var / mean of inning"
$ws.Range("B7").Value = 44645
$ws.Range("C7").Value = "Thomas"
$ws.Range("E7").Value = 44645
$ws.Range("F7").Value = "Thomas_workspace"

# Column B ("Date added") gets a new date number format (m/d/yyyy, numFmtId 14)
# instead of the old d-mmm (numFmtId 16) one. Apply it once, then use
# copy/paste-special (format only) so every cell shares the same style index.
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B3:B7").PasteSpecial(-4122)

# E7 (the relocated "Commited on" value) keeps the original d-mmm format.
$ws.Range("E7").NumberFormat = "d-mmm"

$ws.Range("E38").Select()
